$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values per row (columns C-I, O), case with 380 kV
$data = @{
    2 = @{ "C"=2.834502210481434; "D"=2.921953217054471; "E"=16.70273952699619; "F"=20.05797417919059; "G"=24.47440444073666; "H"=11.41363832979199; "I"=16.38611221003454; "O"=17.07954586972049 }
    3 = @{ "C"=2.808503672603291; "D"=2.921686010819879; "E"=15.73979365400252; "F"=19.55191849635494; "G"=23.37815233946197; "H"=11.34212789423988; "I"=16.0983595678437; "O"=16.76986830385365 }
    4 = @{ "C"=2.7931115702881; "D"=2.921759718578869; "E"=15.12273035622041; "F"=19.2410383337258; "G"=22.6858624696143; "H"=11.30089585243095; "I"=15.9241860918019; "O"=16.58213466568272 }
    5 = @{ "C"=2.786989497204141; "D"=2.921849622287121; "E"=14.86504626354034; "F"=19.11452490033622; "G"=22.39939137988487; "H"=11.2847801259625; "I"=15.85394126840655; "O"=16.50634811581733 }
    6 = @{ "C"=2.785982202230645; "D"=2.921868168734099; "E"=14.82189077943927; "F"=19.09353423208212; "G"=22.35157473898734; "H"=11.28214597632841; "I"=15.84232435492566; "O"=16.49381030925658 }
    7 = @{ "C"=2.793028388484828; "D"=2.921760688527621; "E"=15.11927995754612; "F"=19.23933114346435; "G"=22.68201600144656; "H"=11.30067571247398; "I"=15.92323564679228; "O"=16.58110953531927 }
    8 = @{ "C"=2.825423502009746; "D"=2.92181178128794; "E"=16.37621570653947; "F"=19.88364487339597; "G"=24.10066096499336; "H"=11.38843486064294; "I"=16.28643572825291; "O"=16.97233525539113 }
    9 = @{ "C"=2.893173323068878; "D"=2.923795484154855; "E"=18.77940941808601; "F"=21.13692313316225; "G"=26.71277809985449; "H"=11.58109830804086; "I"=17.01409898298126; "O"=17.75380677738752 }
    10 = @{ "C"=2.945095738912458; "D"=2.926397385555782; "E"=20.45630198673855; "F"=22.04018808408803; "G"=28.50850979450649; "H"=11.73422960676114; "I"=17.55220620562617; "O"=18.33027937365811 }
    11 = @{ "C"=2.969088735992313; "D"=2.927828707778979; "E"=21.17649355917161; "F"=22.44514234115035; "G"=29.29539143617668; "H"=11.80618519361414; "I"=17.79655404344216; "O"=18.59173846750069 }
    12 = @{ "C"=2.978219953219948; "D"=2.928406236134329; "E"=21.44311882955934; "F"=22.59747003669439; "G"=29.58883096260398; "H"=11.83374313777814; "I"=17.88892260391505; "O"=18.69053090742767 }
    13 = @{ "C"=2.976251483190407; "D"=2.928280277094802; "E"=21.38596685349935; "F"=22.56471137841231; "G"=29.52583837487109; "H"=11.82779458430778; "I"=17.86903803392806; "O"=18.66926542142356 }
    14 = @{ "C"=2.969839099306019; "D"=2.927875510157183; "E"=21.198550945294; "F"=22.45769576608214; "G"=29.3196249215489; "H"=11.80844630836659; "I"=17.80415702362779; "O"=18.59987109836102 }
    15 = @{ "C"=2.965917026405712; "D"=2.927632201063894; "E"=21.08296049945112; "F"=22.39200815117388; "G"=29.19271661739474; "H"=11.79663468415417; "I"=17.76439178787495; "O"=18.55733387493649 }
    16 = @{ "C"=2.943534559880374; "D"=2.926308821137681; "E"=20.40838086130207; "F"=22.0135893548844; "G"=28.45646176678647; "H"=11.72957162029899; "I"=17.53622022891922; "O"=18.3131676816905 }
    17 = @{ "C"=2.929893489230276; "D"=2.925560342612827; "E"=19.98365333694488; "F"=21.77979758343483; "G"=27.99693971637658; "H"=11.68900436561274; "I"=17.39606113090698; "O"=18.16310401839636 }
    18 = @{ "C"=2.922083121824703; "D"=2.925153159089977; "E"=19.73534616995748; "F"=21.64477615587814; "G"=27.72982129512861; "H"=11.66588845484039; "I"=17.31540951244098; "O"=18.07672379090912 }
    19 = @{ "C"=2.91944502870715; "D"=2.925019302266265; "E"=19.65058295046881; "F"=21.59897088845626; "G"=27.63890367331788; "H"=11.6580997302975; "I"=17.28809914378474; "O"=18.04746857153799 }
    20 = @{ "C"=2.931341976901226; "D"=2.925637606425362; "E"=20.02928160063101; "F"=21.80474333456367; "G"=28.04614951869005; "H"=11.69330047348332; "I"=17.41098573139984; "O"=18.17908631470509 }
    21 = @{ "C"=2.971721398969819; "D"=2.927993436917682; "E"=21.2537646947785; "F"=22.48915776837043; "G"=29.38031949041873; "H"=11.81412111816856; "I"=17.82321926515934; "O"=18.62026057364823 }
    22 = @{ "C"=2.99837368311897; "D"=2.929740081761629; "E"=22.01852129484445; "F"=22.93045133276988; "G"=30.2257616176663; "H"=11.89488114835647; "I"=18.09165908564866; "O"=18.90728699567929 }
    23 = @{ "C"=2.984127550908814; "D"=2.928788961441041; "E"=21.61359423408222; "F"=22.69552572357222; "G"=29.77702354447645; "H"=11.85162039647263; "I"=17.94850833355133; "O"=18.75424828378111 }
    24 = @{ "C"=2.930687015147872; "D"=2.925602603432695; "E"=20.00866591621071; "F"=21.7934672522137; "G"=28.02391088581592; "H"=11.69135755756652; "I"=17.40423853716604; "O"=18.17186104263799 }
    25 = @{ "C"=2.874436613818316; "D"=2.923057563585532; "E"=18.12383623505978; "F"=20.80018057611107; "G"=26.02655440719849; "H"=11.5268708537711; "I"=16.81621279524715; "O"=17.54153717719697 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 3).Value = $row["C"]
    $ws.Cells.Item($r, 4).Value = $row["D"]
    $ws.Cells.Item($r, 5).Value = $row["E"]
    $ws.Cells.Item($r, 6).Value = $row["F"]
    $ws.Cells.Item($r, 7).Value = $row["G"]
    $ws.Cells.Item($r, 8).Value = $row["H"]
    $ws.Cells.Item($r, 9).Value = $row["I"]
    $ws.Cells.Item($r, 15).Value = $row["O"]
}
